$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38; this shifts the existing rows 38-74
# down to 39-75 (and all their formatting/styles along with them).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly price record.
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = 'Vega Modelo de Temuco'
$ws.Range("C38").Value = 'La Araucanía'
$ws.Range("D38").Value = 44781
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 300000001
$ws.Range("G38").Value = 'Rabanito'
$ws.Range("H38").Value = 'Sin especificar'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 20
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 10000
$ws.Range("N38").Value = '$/docena de paquetes'
$ws.Range("O38").Value = 'Provincia de Cautín'
$ws.Range("P38").Value = 833
$ws.Range("Q38").Value = 12
$ws.Range("R38").Value = 'Hortaliza'
